# The deck's slide-master theme ("Integral", ppt/theme/theme1.xml) is being
# swapped for the stock PowerPoint "Office Theme" palette that previously
# only lived on the Notes Master's theme part (ppt/theme/theme2.xml).
#
# The two themes already share an identical font scheme (fontScheme) and
# format scheme (fmtScheme) — only the 12-slot colour scheme (clrScheme)
# differs between "Integral" and "Office". So re-pointing the slide
# master's theme colours at the "Office Theme" palette reproduces the
# substantive part of the change.
#
# PpColorSchemeIndex-style helper: PowerPoint COM `RGBColor.RGB` is stored
# little-endian (BGR) the same way VBA's RGB() function packs it.
function ConvertTo-ComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

# Order matches the DrawingML <a:clrScheme> child order / ThemeColorScheme
# item order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ConvertTo-ComRgb $officeTheme[$i - 1]
}
